# DU FBS Mock 2 - results update
# - Activates Sheet3 (was Sheet1) as the selected tab / window view
# - Sets the bottom-right pane selection on Sheet3 to CD16
# - Fills in the previously-blank answer cells in row 22 (student "Hridoy
#   Ahmed Nijhum" / similar) for columns L:AT and BK:BZ with graded
#   answer-key strings, reusing the existing "correct" (green) / "wrong"
#   (pink) cell styles already used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Style source cells already present on Sheet3 - C22 carries the "correct"
# (green fill) style, I21 carries the "wrong" (pink fill) style. Copying
# from these reuses the workbook's existing cellXfs entries (s="16" / s="17")
# instead of minting new, duplicate styles.
$correctStyleSrc = $ws3.Range("C22")
$wrongStyleSrc = $ws3.Range("I21")

$cellData = @(
    @{ Cell = "L22"; Style = 16; Value = "D (C)" }
    @{ Cell = "M22"; Style = 16; Value = "B (C)" }
    @{ Cell = "N22"; Style = 16; Value = "C (C)" }
    @{ Cell = "O22"; Style = 16; Value = "D (C)" }
    @{ Cell = "P22"; Style = 16; Value = "C (C)" }
    @{ Cell = "Q22"; Style = 16; Value = "D (C)" }
    @{ Cell = "R22"; Style = 16; Value = "B (C)" }
    @{ Cell = "S22"; Style = 16; Value = "C (C)" }
    @{ Cell = "T22"; Style = 16; Value = "C (C)" }
    @{ Cell = "U22"; Style = 16; Value = "B (C)" }
    @{ Cell = "V22"; Style = 16; Value = "C (C)" }
    @{ Cell = "W22"; Style = 16; Value = "C (C)" }
    @{ Cell = "X22"; Style = 16; Value = "A (C)" }
    @{ Cell = "Y22"; Style = 16; Value = "B (C)" }
    @{ Cell = "Z22"; Style = 16; Value = "A (C)" }
    @{ Cell = "AA22"; Style = 17; Value = "C (W)" }
    @{ Cell = "AB22"; Style = 16; Value = "B (C)" }
    @{ Cell = "AC22"; Style = 17; Value = "A (W)" }
    @{ Cell = "AD22"; Style = 16; Value = "B (C)" }
    @{ Cell = "AE22"; Style = 16; Value = "C (C)" }
    @{ Cell = "AF22"; Style = 16; Value = "C (C)" }
    @{ Cell = "AG22"; Style = 16; Value = "C (C)" }
    @{ Cell = "AH22"; Style = 17; Value = "B (W)" }
    @{ Cell = "AI22"; Style = 17; Value = "D (W)" }
    @{ Cell = "AJ22"; Style = 17; Value = "B (W)" }
    @{ Cell = "AK22"; Style = 17; Value = "B (W)" }
    @{ Cell = "AL22"; Style = 16; Value = "B (C)" }
    @{ Cell = "AM22"; Style = 16; Value = "D (C)" }
    @{ Cell = "AN22"; Style = 17; Value = "A (W)" }
    @{ Cell = "AO22"; Style = 16; Value = "A (C)" }
    @{ Cell = "AP22"; Style = 17; Value = "A (W)" }
    @{ Cell = "AQ22"; Style = 17; Value = "A (W)" }
    @{ Cell = "AR22"; Style = 16; Value = "B (C)" }
    @{ Cell = "AS22"; Style = 16; Value = "A (C)" }
    @{ Cell = "AT22"; Style = 17; Value = "B (W)" }
    @{ Cell = "BK22"; Style = 17; Value = "C (W)" }
    @{ Cell = "BL22"; Style = 17; Value = "B (W)" }
    @{ Cell = "BM22"; Style = 16; Value = "C (C)" }
    @{ Cell = "BN22"; Style = 17; Value = "C (W)" }
    @{ Cell = "BO22"; Style = 16; Value = "C (C)" }
    @{ Cell = "BP22"; Style = 17; Value = "D (W)" }
    @{ Cell = "BQ22"; Style = 16; Value = "C (C)" }
    @{ Cell = "BR22"; Style = 16; Value = "D (C)" }
    @{ Cell = "BS22"; Style = 16; Value = "B (C)" }
    @{ Cell = "BT22"; Style = 17; Value = "B (W)" }
    @{ Cell = "BU22"; Style = 16; Value = "B (C)" }
    @{ Cell = "BW22"; Style = 17; Value = "C (W)" }
    @{ Cell = "BX22"; Style = 16; Value = "C (C)" }
    @{ Cell = "BY22"; Style = 16; Value = "C (C)" }
    @{ Cell = "BZ22"; Style = 16; Value = "B (C)" }
)

foreach ($entry in $cellData) {
    $target = $ws3.Range($entry.Cell)
    if ($entry.Style -eq 16) {
        $correctStyleSrc.Copy($target)
    } else {
        $wrongStyleSrc.Copy($target)
    }
    $target.Value = $entry.Value
}

# Switch the active/selected tab from Sheet1 to Sheet3 and set the new
# selection on Sheet3's frozen-pane view (bottom-right pane -> CD16).
$ws3.Activate()
$ws3.Range("CD16").Select()
